# Re-shuffle the species-observation rows 69-77,79,81-82 on the "Artfynd"
# sheet. Columns A (Id), B (Taxonsorteringsordning), D (Rodlistade),
# E (TaxonId), F (Artnamn), G (Vetenskapligt namn), H (Auktor), Q (Ost) and
# R (Nord) get redistributed across those rows while every other column
# (locality, dates, observer, ...) stays put on its original row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per target row, taken straight from the target workbook state.
$data = @{
    69 = @{ A = 111785199; B = 89416; D = "LC"; E = 1205;   F = "Stor aspticka"; G = "Phellinus populicola";           H = "Niemelä";                       Q = 577256.110519147;  R = 6944531.123615563 }
    70 = @{ A = 111785229; B = 78578; D = "NT"; E = 6458;   F = "Lunglav";        G = "Lobaria pulmonaria";            H = "(L.) Hoffm.";                   Q = 577208.3826684169; R = 6944521.722980071 }
    71 = @{ A = 111785201; B = 78512; D = "LC"; E = 6456;   F = "Skinnlav";       G = "Leptogium saturninum";          H = "(Dicks.) Nyl.";                 Q = 577248.2772659193; R = 6944530.940753835 }
    72 = @{ A = 111785206; B = 77268; D = "NT"; E = 228912; F = "Mörk kolflarnlav"; G = "Carbonicola myrmecina";       H = "(Ach.) Bendiksby & Timdal";     Q = 577235.6798241453; R = 6944655.86623876 }
    73 = @{ A = 111785192; B = 89405; D = "NT"; E = 1202;   F = "Ullticka";       G = "Phellinidium ferrugineofuscum"; H = "(P.Karst.) Fiasson & Niemelä";  Q = 577281.7951240344; R = 6944714.487089146 }
    74 = @{ A = 111785200; B = 78512; D = "LC"; E = 6456;   F = "Skinnlav";       G = "Leptogium saturninum";          H = "(Dicks.) Nyl.";                 Q = 577256.110519147;  R = 6944531.123615563 }
    75 = @{ A = 111785244; B = 96348; D = "VU"; E = 220787; F = "Knärot";         G = "Goodyera repens";               H = "(L.) R. Br.";                   Q = 577364.1822193989; R = 6944621.578847388 }
    76 = @{ A = 111785251; B = 93161; D = "VU"; E = 1079;   F = "Aspfjädermossa"; G = "Neckera pennata";               H = "Hedw.";                         Q = 577283.2535308318; R = 6944533.598891968 }
    77 = @{ A = 111785230; B = 78578; D = "NT"; E = 6458;   F = "Lunglav";        G = "Lobaria pulmonaria";            H = "(L.) Hoffm.";                   Q = 577261.8704127767; R = 6944620.109213427 }
    79 = @{ A = 111785202; B = 78512; D = "LC"; E = 6456;   F = "Skinnlav";       G = "Leptogium saturninum";          H = "(Dicks.) Nyl.";                 Q = 577215.0430418774; R = 6944631.445974576 }
    81 = @{ A = 111785228; B = 78578; D = "NT"; E = 6458;   F = "Lunglav";        G = "Lobaria pulmonaria";            H = "(L.) Hoffm.";                   Q = 577256.110519147;  R = 6944531.123615563 }
    82 = @{ A = 111785235; B = 77267; D = "NT"; E = 6446;   F = "Kolflarnlav";    G = "Carbonicola anthracophila";     H = "(Nyl.) Bendiksby & Timdal";     Q = 577226.625646919;  R = 6944648.749557905 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 1).Value  = $vals.A   # A: Id
    $ws.Cells.Item($row, 2).Value  = $vals.B   # B: Taxonsorteringsordning
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D: Rodlistade
    $ws.Cells.Item($row, 5).Value  = $vals.E   # E: TaxonId
    $ws.Cells.Item($row, 6).Value  = $vals.F   # F: Artnamn
    $ws.Cells.Item($row, 7).Value  = $vals.G   # G: Vetenskapligt namn
    $ws.Cells.Item($row, 8).Value  = $vals.H   # H: Auktor
    $ws.Cells.Item($row, 17).Value = $vals.Q   # Q: Ost
    $ws.Cells.Item($row, 18).Value = $vals.R   # R: Nord
}
